# Updates cryptos list (price / 1h volume columns) to reflect the
# latest scrape, as produced by the "Updated cryptos list ... with GitHub
# Actions" commit. Also fixes the HuobiToken/MXToken row ordering (rows 36
# and 37 swap places in the source data).
#
# D-column prices that look like plain numbers (e.g. "284.13") would be
# auto-coerced to a numeric type by a bare Range.Value assignment, which
# would not match the source workbook's inline-string cells. Forcing the
# NumberFormat to text ("@") immediately before the assignment keeps the
# write as text, and ClearFormats() afterwards drops the temporary format
# again so no stray style is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.897.28'
$ws.Range("E2").Value = '  +5.31%  '
$ws.Range("D3").Value = '1.878.86'
$ws.Range("E3").Value = '  +4.04%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '284.13'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.00%  '
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5222'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +4.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3548'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.33'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +3.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07089'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +6.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.48'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.8259'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.99%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07757'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.57%  '
$ws.Range("D14").Value = '1.875.65'
$ws.Range("E14").Value = '  +3.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.216'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +3.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.30'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +3.55%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9985'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.54'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +4.64%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008185'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +3.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9986'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.13%  '
$ws.Range("D21").Value = '26.930.15'
$ws.Range("E21").Value = '  +5.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.796'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.86%  '
$ws.Range("E23").Value = '  +2.69%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.263'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +3.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.435'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +15.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '145.76'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +3.05%  '
$ws.Range("E27").Value = '  +3.39%  '
$ws.Range("E28").Value = '  +0.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '112.06'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +3.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.434'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +3.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.392'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +4.46%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08889'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04949'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.91%  '
$ws.Range("E34").Value = '  +5.68%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7551'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +3.54%  '
$ws.Range("B36").Value = 'MXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.299'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +9.27%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.869'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.426'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +6.67%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5368'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +3.78%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01892'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.95%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9872'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.86%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '117.19'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +3.76%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.332'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.272'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.91%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4673'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9981'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.25%  '
$ws.Range("E47").Value = '  -0.31%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.508'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +3.45%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.91'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +3.50%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.526'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.60%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05945'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.37%  '
